# Apply the commit's changes to config.xlsx
#
# Summary of the edit:
#  - On the "Feeds" sheet, row 25 (Slovenia / "SI") loses its feed-URL
#    values in columns B, C and D - only the country code in A25 remains.
#  - The "Feeds" sheet becomes the active sheet/tab (it was "Schedule"
#    before), and its remembered selection moves from C38 to B25.

$wb = $excel.ActiveWorkbook
$feeds = $wb.Worksheets.Item("Feeds")

# Remove the Slovenian feed URLs (B25:D25); keep the country code in A25.
$feeds.Range("B25:D25").ClearContents() | Out-Null

# Update the remembered selection on the Feeds sheet and make it active
# (previously "Schedule" was the active/selected tab).
$feeds.Range("B25").Select() | Out-Null
$feeds.Activate() | Out-Null
